# A new weekly price observation (2021-10-25, "Región de Ñuble") is inserted
# as a new data row right before the current row 88. Every existing record
# from row 88 down to the last row (191) shifts down by one row, and the
# dimension grows from A1:R191 to A1:R192.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 88 — pushes rows 88:191 down to 89:192
# and copies formatting (incl. the date number format on column D) from the
# row above, same as Excel's default Insert behaviour.
$ws.Rows.Item(88).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(88, 1).Value  = 7
$ws.Cells.Item(88, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(88, 3).Value  = "Ñuble"
$ws.Cells.Item(88, 4).Value  = 44494
$ws.Cells.Item(88, 5).Value  = 16
$ws.Cells.Item(88, 6).Value  = 100114013
$ws.Cells.Item(88, 7).Value  = "Zanahoria"
$ws.Cells.Item(88, 8).Value  = "Sin especificar"
$ws.Cells.Item(88, 9).Value  = "Primera"
$ws.Cells.Item(88, 10).Value = 120
$ws.Cells.Item(88, 11).Value = 8500
$ws.Cells.Item(88, 12).Value = 9000
$ws.Cells.Item(88, 13).Value = 8750
$ws.Cells.Item(88, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(88, 15).Value = "Región de Ñuble"
$ws.Cells.Item(88, 16).Value = 438
$ws.Cells.Item(88, 17).Value = 20
$ws.Cells.Item(88, 18).Value = "Hortaliza"
